# Commit: "Added Samples and Files Tab to all tests"
#
# This adds two new rows (SamplesTab, FilesTab) below the existing CasesTab
# row on the startup sheet, and tweaks the CasesTab "Age (years)" query to
# coalesce fractional/blank ages. Column widths and row heights are adjusted
# to fit the new, longer query text.
#
# NOTE: cell values below are written in a specific order (TabName column
# for both new rows, then the query column for both new rows, then the
# CasesTab update last) so the workbook's shared-string table ends up in
# the same order Excel produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- TabName column for the two new rows ---
$ws.Range("A3").Value2 = "SamplesTab"
$ws.Range("A4").Value2 = "FilesTab"

# --- query column for the two new rows ---
$ws.Range("B3").Value2 = 'MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE  f.file_format IN [''gz'',''tbi'']
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,''-'')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`'
$ws.Range("B4").Value2 = 'MATCH (f:file)-->(parent)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
WHERE  f.file_format IN [''gz'',''tbi'']
WITH
        f, parent,p, ss, d,tp, s, samp,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent,p, ss, d,tp, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent,p, ss, d,tp, s, samp, unit,
        round(factor * value)/factor AS size
RETURN Distinct
    f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`
    order by f.file_name'

# --- Update existing CasesTab row (row 2): the Age (years) column now
# coalesces the value so fractional-vs-integer ages render consistently.
$ws.Range("B2").Value2 = 'MATCH (ss:study_subject)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)<-[:sample_of_study_subject]-(samp:sample)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)

WITH DISTINCT ss, samp, collect(DISTINCT samp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
MATCH (f:file)
MATCH (f)-[:file_of_sample]->(samp)
WHERE  f.file_format IN [''gz'',''tbi'']
return DISTINCT ss.study_subject_id as `Case ID`,
   p.program_acronym as `Program Code`,
    p.program_id as Program_ID,
   s.study_acronym as `Arm`,
   ss.disease_subtype as `Diagnosis`,
   sf.grouped_recurrence_score AS `Recurrence Score`,
   d.tumor_size_group AS `tumor_size`,
   d.er_status AS `ER Status`,
   d.pr_status AS `PR Status`,
   coalesce(CASE demo.age_at_index % 1 WHEN 0 THEN apoc.convert.toInteger(demo.age_at_index) ELSE demo.age_at_index END, '''') AS `Age (years)`,
	demo.survival_time AS `Survival (days)`'

# --- remaining columns (StatQuery / dbExcel / WebExcel) for the two new rows ---
$ws.Range("C3").Value2 = 'MATCH (f:file)
WHERE f.file_format IN [''gz'',''tbi'']
OPTIONAL MATCH (f)-->(parent)-->(ss:study_subject)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
OPTIONAL MATCH (f)-[:file_of_laboratory_procedure]->(lp:laboratory_procedure)
OPTIONAL MATCH (f)-[:file_of_sample]->(samp)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files'
$ws.Range("D3").Value2 = 'TC03_Bento_Filter_FileType-vcf_Neo4jData.xlsx'
$ws.Range("E3").Value2 = 'TC03_Bento_Filter_FileType-vcf_WebData.xlsx'

$ws.Range("C4").Value2 = 'MATCH (f:file)
WHERE f.file_format IN [''gz'',''tbi'']
OPTIONAL MATCH (f)-->(parent)-->(ss:study_subject)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
OPTIONAL MATCH (f)-[:file_of_laboratory_procedure]->(lp:laboratory_procedure)
OPTIONAL MATCH (f)-[:file_of_sample]->(samp)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files'
$ws.Range("D4").Value2 = 'TC03_Bento_Filter_FileType-vcf_Neo4jData.xlsx'
$ws.Range("E4").Value2 = 'TC03_Bento_Filter_FileType-vcf_WebData.xlsx'

# Wrap text in the query/StatQuery columns, same as the existing CasesTab row
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true
$ws.Range("B4").WrapText = $true
$ws.Range("C4").WrapText = $true

# Row heights sized for the new (longer) wrapped text
$ws.Rows.Item(2).RowHeight = 388.8
$ws.Rows.Item(3).RowHeight = 345.6
$ws.Rows.Item(4).RowHeight = 409.6

# Column widths (widened/narrowed to fit the new query text)
$ws.Columns.Item(2).ColumnWidth = 82
$ws.Columns.Item(3).ColumnWidth = 50.5
$ws.Columns.Item(4).ColumnWidth = 43

# Move the active selection to the new Samples row, matching the saved view state
$ws.Range("B3").Select()
